$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.869.12'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '2.383.81'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '555.22'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '133.26'
$ws.Range('E6').Value = '  -2.83%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '5.63'
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('E12').Value = '  -2.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '24.41'
$ws.Range('E13').Value = '  -4.51%  '
$ws.Range('D14').Value = '2.809.19'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').Value = '59.816.59'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '2.380.90'
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '4.48'
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '320.35'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '64.19'
$ws.Range('E23').Value = '  -3.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '0.172'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '8.41'
$ws.Range('E26').Value = '  -2.53%  '
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('E29').Value = '  -1.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '169.67'
$ws.Range('E30').Value = '  +0.84%  '
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('E32').Value = '  +5.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '0.397'
$ws.Range('E33').Value = '  -2.87%  '
$ws.Range('E34').Value = '  -2.39%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '318.65'
$ws.Range('E39').Value = '  +1.61%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '1.57'
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '38.63'
$ws.Range('E41').Value = '  -2.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '145.98'
$ws.Range('E42').Value = '  +5.32%  '
$ws.Range('E43').Value = '  -4.21%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '19.68'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('E48').Value = '  -2.91%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '0.947'
$ws.Range('E51').Value = '  -0.19%  '
